$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original data block (product rows) occupies A697:T795 (99 rows).
# A new weekly entry (3 rows: Especial/Primera/Segunda) is inserted at row 697,
# pushing the existing 99 rows down by 3 (to rows 700:798).

$srcRange = "A697:T795"
$destRange = "A700:T798"

# 1) Read the existing 99-row block before we overwrite anything.
$block = $ws.Range($srcRange).Value()

# 2) Preserve the date-column (D) number format so it carries over to the
#    freshly-created rows 796:798 (otherwise Excel would default to a plain
#    date format for brand new cells). Apply the format before writing the
#    values so no transient/auto style gets created along the way.
$dateFormat = $ws.Range("D697").NumberFormat
$ws.Range("D700:D798").NumberFormat = $dateFormat

# 3) Shift the whole block down by 3 rows.
$ws.Range($destRange).Value = $block

$ws.Range("D697").Value = 44748
$ws.Range("N697").Value = 7500
$ws.Range("O697").Value = 8000
$ws.Range("P697").Value = 7750
$ws.Range("R697").Value = "Brasil"
$ws.Range("S697").Value = 1938

$ws.Range("D698").Value = 44748
$ws.Range("N698").Value = 7500
$ws.Range("O698").Value = 8000
$ws.Range("P698").Value = 7750
$ws.Range("R698").Value = "Brasil"
$ws.Range("S698").Value = 1938

$ws.Range("D699").Value = 44748
$ws.Range("N699").Value = 7500
$ws.Range("O699").Value = 8000
$ws.Range("P699").Value = 7750
$ws.Range("R699").Value = "Brasil"
$ws.Range("S699").Value = 1938
